$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New historical date rows (column A only), continuing the existing series.
# Use Copy/PasteSpecial(formats) from the last existing date cell (A370) so the
# new cells inherit the same date style (s="2") already used by the column,
# then fill in the actual date serial values.
$ws.Range("A370").Copy()
$ws.Range("A371:A378").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$dates = @(45492, 45493, 45494, 45495, 45496, 45497, 45498, 45499)
for ($i = 0; $i -lt $dates.Length; $i++) {
    $ws.Cells.Item(371 + $i, 1).Value = $dates[$i]
}

# Match the new selection left by the edit (activeCell A373, multi-row sqref).
$ws.Range("A373:A378").Select()
